# cond01.xlsx update:
#  - reduce stimulus size (sizeL/sizeR, cols F/G) from 6.7 to 4 deg
#  - reduce spatial frequency (sfL/sfR, cols H/I) from 1.5 to 1 cyc/deg
#  - for all 16 trial rows (rows 2-17)
#  - update window zoom and active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 6).Value = 4   # F: sizeL
    $ws.Cells.Item($r, 7).Value = 4   # G: sizeR
    $ws.Cells.Item($r, 8).Value = 1   # H: sfL
    $ws.Cells.Item($r, 9).Value = 1   # I: sfR
}

# Zoom in to 145% and move the active selection to K5
$excel.ActiveWindow.Zoom = 145
$ws.Range("K5").Select()
